$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The ranking table got recomputed: rows 4/5 swap identity (Matthew <-> Drew) and
# rows 6/7 swap identity (Juan <-> Edosagbe), each person keeping their own
# prolific number (B) and race (G) as they move row, while every realeffort
# score (F) is refreshed with a newly computed value.

$ws.Cells.Item(2, 3).Value  = "60b091ed11ccda59e3fc7761"
$ws.Cells.Item(2, 4).Value  = "Myles"
$ws.Cells.Item(2, 2).Value  = 44
$ws.Cells.Item(2, 6).Value  = 13.19982871425305
$ws.Cells.Item(2, 7).Value  = "Black or African American"

$ws.Cells.Item(3, 3).Value  = "601d69a993d94008fb2b25dc"
$ws.Cells.Item(3, 4).Value  = "Quinterius"
$ws.Cells.Item(3, 2).Value  = 3
$ws.Cells.Item(3, 6).Value  = 8.081433205567341
$ws.Cells.Item(3, 7).Value  = "Black or African American"

$ws.Cells.Item(4, 3).Value  = "60c2341fe95d71ee52c043f0"
$ws.Cells.Item(4, 4).Value  = "Matthew"
$ws.Cells.Item(4, 2).Value  = 30
$ws.Cells.Item(4, 6).Value  = 7.048241805590385
$ws.Cells.Item(4, 7).Value  = "White"

$ws.Cells.Item(5, 3).Value  = "5ff8ad350d084e10f500e48a"
$ws.Cells.Item(5, 4).Value  = "Drew"
$ws.Cells.Item(5, 2).Value  = 27
$ws.Cells.Item(5, 6).Value  = 7.000947600168775
$ws.Cells.Item(5, 7).Value  = "White"

$ws.Cells.Item(6, 3).Value  = "5dd671942b033b5ec8bc97b4"
$ws.Cells.Item(6, 4).Value  = "Juan"
$ws.Cells.Item(6, 2).Value  = 26
$ws.Cells.Item(6, 6).Value  = 5.239313832273305
$ws.Cells.Item(6, 7).Value  = "Hispanic"

$ws.Cells.Item(7, 3).Value  = "60db4fde6193c50664c9c478"
$ws.Cells.Item(7, 4).Value  = "Edosagbe"
$ws.Cells.Item(7, 2).Value  = 22
$ws.Cells.Item(7, 6).Value  = 5.186302527479196
$ws.Cells.Item(7, 7).Value  = "Black or African American"

$ws.Cells.Item(8, 3).Value  = "60bf9943e4e04642d4634ecc"
$ws.Cells.Item(8, 4).Value  = "Jamarii"
$ws.Cells.Item(8, 2).Value  = 32
$ws.Cells.Item(8, 6).Value  = 5.018452747422359
$ws.Cells.Item(8, 7).Value  = "Black or African American"

$ws.Cells.Item(9, 3).Value  = "5e2522d6b734b47915f88275"
$ws.Cells.Item(9, 4).Value  = "Corey"
$ws.Cells.Item(9, 2).Value  = 2
$ws.Cells.Item(9, 6).Value  = 4.155549573790759
$ws.Cells.Item(9, 7).Value  = "White"

$ws.Cells.Item(10, 3).Value = "60b322994d0b901954690036"
$ws.Cells.Item(10, 4).Value = "Brennan"
$ws.Cells.Item(10, 2).Value = 33
$ws.Cells.Item(10, 6).Value = 4.034175108618071
$ws.Cells.Item(10, 7).Value = "White"

$ws.Cells.Item(11, 3).Value = "6088fc724afd5c008db33e9d"
$ws.Cells.Item(11, 4).Value = "Masuf"
$ws.Cells.Item(11, 2).Value = 49
$ws.Cells.Item(11, 6).Value = 3.262000135003892
$ws.Cells.Item(11, 7).Value = "Asian"

$ws.Cells.Item(12, 3).Value = "60b83826821417f8e484a207"
$ws.Cells.Item(12, 4).Value = "Eli"
$ws.Cells.Item(12, 2).Value = 29
$ws.Cells.Item(12, 6).Value = 2.30063950806506
$ws.Cells.Item(12, 7).Value = "White"

$ws.Cells.Item(13, 3).Value = "6097b95056caf5ebb2720002"
$ws.Cells.Item(13, 4).Value = "Damian"
$ws.Cells.Item(13, 2).Value = 50
$ws.Cells.Item(13, 6).Value = 2.244526961475056
$ws.Cells.Item(13, 7).Value = "Black or African American"
